$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# A leading apostrophe forces Excel to store purely-numeric-looking
# strings (e.g. '529.77') as text instead of coercing them to a float,
# and resetting .Style to 'Normal' afterwards clears the auto-applied
# quote-prefix/text style so the cell keeps its original (default) style.

$ws.Range('D2').Value = "'69.765.23"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.38%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.942.18"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.93%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.08%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'529.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'146.63"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.15%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.28%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.06%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.731"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.57%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +5.26%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0000343"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.25%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'42.91"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.15%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'10.55"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.78%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'4.579.35"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.34%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.941.85"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.00%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'14.18"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.88%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.16%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  +7.09%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'19.89"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.56%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'69.734.42"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +2.22%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'435.11"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.45%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -4.11%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -2.62%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'88.59"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.51%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = "'PancakeSwap"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'4.05"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +11.92%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = "'RenderToken"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'11.95"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +5.46%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'11.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.40%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'36.77"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -3.59%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.08%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'703.68"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.72%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'13.36"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -2.94%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.128"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -2.11%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.20%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'67.92"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +12.30%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.443"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +9.27%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value = "'NEARProtocol"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'6.07"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -2.95%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = "'PEPE"
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'0.0₃0874"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.95%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'40.45"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -3.09%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.45%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.999"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.03%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.07%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.31%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'2.85"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -3.63%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'3.09"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +6.41%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -4.04%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'3.26"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +15.93%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'3.42"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +2.59%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +1.50%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.0₆0366"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +7.56%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.17%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.71%  "
$ws.Range('E51').Style = 'Normal'
